# Aggiornamento dati fino al 9/09 compreso: aggiunge le righe 367-374
# (giorni 2021-09-02 .. 2021-09-09) in fondo al foglio, replicando lo
# stile della colonna A (formato data) gia' usato per le righe precedenti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 366

$newData = @(
    @(44441, 1, 4, 75.75757575757575),
    @(44442, 1, 4, 75.75757575757575),
    @(44443, 0, 3, 56.81818181818181),
    @(44444, 0, 3, 56.81818181818181),
    @(44445, 0, 2, 37.87878787878788),
    @(44446, 0, 2, 37.87878787878788),
    @(44447, 0, 2, 37.87878787878788),
    @(44448, 0, 1, 18.93939393939394)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $targetRow = $lastRow + 1 + $i
    $srcRow = $targetRow - 1

    # Copia il formato della colonna A (data, bordo, grassetto, ecc.)
    # dalla riga precedente, cosi' la nuova riga eredita lo stile "s=2".
    $ws.Range("A$srcRow").Copy() | Out-Null
    $ws.Range("A$targetRow").PasteSpecial(-4122) | Out-Null

    $row = $newData[$i]
    $ws.Cells.Item($targetRow, 1).Value = $row[0]
    $ws.Cells.Item($targetRow, 2).Value = $row[1]
    $ws.Cells.Item($targetRow, 3).Value = $row[2]
    $ws.Cells.Item($targetRow, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
